$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in O1: "动态寻址（可选）" -> "寻址方式（可选）"
$ws.Range("O1").Value = "寻址方式（可选）"

# O2 used to be a boolean TRUE switch; now becomes text "静态"
$ws.Range("O2").Value = "静态"

# O3 is a new cell with text "动态"
$ws.Range("O3").Value = "动态"

# Update selection to O3 (matches the selection change in the diff)
$ws.Range("O3").Select()

